# Applies the "Updated cryptos list" price/volume refresh (and the
# Algorand/Aave row swap in rows 48-49) described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D/E
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.245.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.46%  "

# Row 3: update D/E
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.605.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +11.35%  "

# Row 4: update D/E
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.53%  "

# Row 5: update D/E
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.37"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.35%  "

# Row 6: update D/E
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.46"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.05%  "

# Row 7: update D/E
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +6.60%  "

# Row 8: update E
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.35%  "

# Row 9: update D/E
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.583"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +14.09%  "

# Row 10: update D/E
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.15"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +12.86%  "

# Row 11: update D/E
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0840"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.15%  "

# Row 12: update D/E
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.00"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +12.56%  "

# Row 13: update D/E
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.011.87"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +11.84%  "

# Row 14: update E
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.86%  "

# Row 15: update D/E
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.690.17"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +15.45%  "

# Row 16: update D/E
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.904"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +12.49%  "

# Row 17: update D/E
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "15.23"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +11.36%  "

# Row 18: update D/E
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.561.95"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.31%  "

# Row 19: update D/E
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.87"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +9.36%  "

# Row 20: update E
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.58%  "

# Row 21: update D/E
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.67"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +11.69%  "

# Row 22: update D/E
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.46"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.55%  "

# Row 23: update D/E
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.34"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.00%  "

# Row 24: update D/E
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.09%  "

# Row 25: update D/E
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.13"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +11.63%  "

# Row 26: update E
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.07%  "

# Row 27: update D/E
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.11"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +25.35%  "

# Row 28: update D/E
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.95"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.60%  "

# Row 29: update D/E
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.56"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +9.16%  "

# Row 30: update D/E
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.59%  "

# Row 31: update D/E
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.74%  "

# Row 32: update D/E
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.01"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.38%  "

# Row 33: update D/E
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.99"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +10.12%  "

# Row 34: update D/E
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.25"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +24.93%  "

# Row 35: update D/E
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0847"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +9.52%  "

# Row 36: update D/E
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "150.31"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.97%  "

# Row 37: update D/E
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.122"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +9.50%  "

# Row 38: update D/E
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.121"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.18%  "

# Row 39: update D/E
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.51"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.47%  "

# Row 40: update D/E
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.29"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +10.92%  "

# Row 41: update D/E
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0333"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +11.71%  "

# Row 42: update D/E
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.65"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +13.97%  "

# Row 43: update D/E
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.025.06"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +9.39%  "

# Row 44: update D/E
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.36%  "

# Row 45: update D/E
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.72"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.76%  "

# Row 46: update D/E
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.97"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +38.08%  "

# Row 47: update D/E
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.90"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.65%  "

# Row 48: update B/C/D/E
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "109.02"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +13.26%  "

# Row 49: update B/C/D/E
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.202"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +8.89%  "

# Row 50: update D/E
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.869.22"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +11.84%  "

# Row 51: update D/E
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.95"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +13.10%  "
